$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.136.99'
$ws.Range('E2').Value = '  +4.82%  '
$ws.Range('D3').Value = '2.674.22'
$ws.Range('E3').Value = '  +8.05%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '114.59'
$ws.Range('E5').Value = '  +9.47%  '
$ws.Range('D6').Value = '327.53'
$ws.Range('E6').Value = '  +3.79%  '
$ws.Range('E7').Value = '  +2.63%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '0.560'
$ws.Range('E9').Value = '  +5.07%  '
$ws.Range('D10').Value = '41.47'
$ws.Range('E10').Value = '  +7.23%  '
$ws.Range('D11').Value = '20.21'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Value = '0.0829'
$ws.Range('E12').Value = '  +3.97%  '
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '7.42'
$ws.Range('E14').Value = '  +5.59%  '
$ws.Range('D15').Value = '3.088.84'
$ws.Range('E15').Value = '  +7.86%  '
$ws.Range('D16').Value = '2.681.26'
$ws.Range('E16').Value = '  +7.19%  '
$ws.Range('D17').Value = '0.881'
$ws.Range('E17').Value = '  +7.13%  '
$ws.Range('D18').Value = '50.048.79'
$ws.Range('E18').Value = '  +4.77%  '
$ws.Range('D19').Value = '13.34'
$ws.Range('E19').Value = '  +5.79%  '
$ws.Range('D20').Value = '6.83'
$ws.Range('E20').Value = '  +4.95%  '
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +4.24%  '
$ws.Range('D23').Value = '281.86'
$ws.Range('E23').Value = '  +3.39%  '
$ws.Range('D24').Value = '72.89'
$ws.Range('E24').Value = '  +3.40%  '
$ws.Range('E25').Value = '  +4.64%  '
$ws.Range('D26').Value = '27.02'
$ws.Range('E26').Value = '  +6.04%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  +4.90%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '36.90'
$ws.Range('E30').Value = '  +7.76%  '
$ws.Range('D31').Value = '0.144'
$ws.Range('E31').Value = '  +4.86%  '
$ws.Range('D32').Value = '50.43'
$ws.Range('E32').Value = '  +2.64%  '
$ws.Range('D33').Value = '5.53'
$ws.Range('E33').Value = '  +5.92%  '
$ws.Range('D34').Value = '19.83'
$ws.Range('E34').Value = '  +5.10%  '
$ws.Range('D35').Value = '0.0819'
$ws.Range('E35').Value = '  +6.83%  '
$ws.Range('E36').Value = '  +13.83%  '
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = '2.09'
$ws.Range('E38').Value = '  +8.83%  '
$ws.Range('E39').Value = '  +10.50%  '
$ws.Range('D40').Value = '125.04'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').Value = '0.114'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '22.45'
$ws.Range('E42').Value = '  +2.19%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.23'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  +7.05%  '
$ws.Range('D45').Value = '2.114.06'
$ws.Range('E45').Value = '  +6.17%  '
$ws.Range('D46').Value = '3.36'
$ws.Range('E46').Value = '  +7.07%  '
$ws.Range('E47').Value = '  +14.88%  '
$ws.Range('E48').Value = '  +5.77%  '
$ws.Range('D49').Value = '9.15'
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('D50').Value = '5.40'
$ws.Range('E50').Value = '  +5.50%  '
$ws.Range('D51').Value = '60.13'
$ws.Range('E51').Value = '  +7.74%  '
